$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update hashcode values (column B) for the rows whose metadata ID (column A) changed.
$updates = @(
    @{ Row = 9; Id = "05-050305TC"; OldHash = "a0c1161837786ed577bd398f7504ad26"; NewHash = "ddb0a2e1826297db71eacc4dbae7632e" }
    @{ Row = 11; Id = "05-050301A"; OldHash = "5e3060384c34d729b00004f7ae145871"; NewHash = "f1b8f9386df95e332d298fbab3b81b20" }
    @{ Row = 15; Id = "05-050207TP"; OldHash = "f70941508953fb002528f6913831f208"; NewHash = "e3d2c480789241c7b9aaad710ce56565" }
    @{ Row = 17; Id = "05-050305TP"; OldHash = "dbb17ca4b52a4c7e5e94472e9b66584d"; NewHash = "599ca733de93e7c04d544d5c4d52a2e4" }
    @{ Row = 24; Id = "05-050316TC"; OldHash = "b5f6d0d26190225d3934df2a1696f2eb"; NewHash = "4dd7784ffc89afbc74a62b1b4e2d44fb" }
    @{ Row = 29; Id = "05-050302A"; OldHash = "40e9a17fbb91f23587413a1d9790a202"; NewHash = "e50a204cc11d6b6acc449796522dee00" }
    @{ Row = 34; Id = "05-050316TP"; OldHash = "e0c01ed14f14658dee48c59beec4108d"; NewHash = "b755446c5707adc7b49b5c01651d7567" }
    @{ Row = 121; Id = "05-050301TP"; OldHash = "89e31980121a03ecb4d90a72f238e8a1"; NewHash = "864991d110708b2f79025d12940c4d2b" }
    @{ Row = 133; Id = "05-050312TP"; OldHash = "d5cb4a103f42f397f487757d5d1944b5"; NewHash = "ef7b25371103721dda82508c20e0d5d1" }
    @{ Row = 136; Id = "05-050312TC"; OldHash = "fdf78d44d331f49c30cdad3e5f603af7"; NewHash = "a561ccab3a1512203e6859ee7abd713e" }
    @{ Row = 159; Id = "05-050203TP"; OldHash = "567cf77756c9ad1d2efe5d8d378938af"; NewHash = "a2346cfada1fded0936ef1f99bc46b31" }
    @{ Row = 162; Id = "05-050308A"; OldHash = "f26197f222ebf1ddb0efdcaf398412ff"; NewHash = "c5578f51422585f435b4795fab343dfb" }
    @{ Row = 169; Id = "05-050203TC"; OldHash = "bd9aecd057b8b0de503941ed6157bbd2"; NewHash = "a4aedcb741238e50c1d8118791741934" }
    @{ Row = 175; Id = "05-050303TP"; OldHash = "2ac35ee76222d5df5e0654457b0c544f"; NewHash = "7eafc91c7445130e9d3c800b2cbaefc2" }
    @{ Row = 180; Id = "05-050303TC"; OldHash = "5d06c54d6ae33de6ef32267596888733"; NewHash = "73775cbc9c547b89b310c1f2030a6f58" }
    @{ Row = 183; Id = "05-050305A"; OldHash = "bfd43c2f789ae217aee9d6a0c58b3db0"; NewHash = "d087296d2235f723f24a0fe0a9d108fe" }
    @{ Row = 191; Id = "05-050314TP"; OldHash = "4d6ab91e2c46180e790ad8b177c98dfb"; NewHash = "3041e26a7928b0249067df8f5ae99ff3" }
    @{ Row = 198; Id = "05-050314TC"; OldHash = "e092507be44a5fc8a934d3c321d61312"; NewHash = "c88e43bb0f2dbec7ce4c2bc1d8f75be3" }
    @{ Row = 200; Id = "05-050306A"; OldHash = "d390d1e05d7bb974a4a42141255c0a4f"; NewHash = "ca1b14fca26ff83e2693d9e224bea198" }
    @{ Row = 213; Id = "05-050303A"; OldHash = "6237ab96bee52eeaaab382477910cf24"; NewHash = "0373eb1d392cfba631626b05d59df8db" }
    @{ Row = 227; Id = "05-050205TP"; OldHash = "38806a07acc53d2c72acefaafded9c2a"; NewHash = "4c77795864a792b7897f0dbb55a7ac1b" }
    @{ Row = 228; Id = "05-050304A"; OldHash = "32971f9e01b7e44aa184d2c517c589a3"; NewHash = "9b5390073171029c530362c11f7e4242" }
    @{ Row = 232; Id = "05-050205TC"; OldHash = "a160247022553052b8abae884d3456af"; NewHash = "4e38ccf629cab9ec45a580a0b4228321" }
    @{ Row = 339; Id = "05-050201TP"; OldHash = "95dbda5d9a8b6ad8dfae2c4599d555fd"; NewHash = "57386898843480ebd5523c50da4cf754" }
    @{ Row = 420; Id = "05-0709-070905BTC"; OldHash = "0841f66eec1f7caf51680bed6f5054c6"; NewHash = "930e9bd628ccd09c643cd2b4a4b8cfad" }
    @{ Row = 464; Id = "05-050204A"; OldHash = "74c3f3300493c90cdb38ddc368e0b62f"; NewHash = "4c4cf74088f38113cd7e5aae873e50f9" }
    @{ Row = 465; Id = "05-050313A"; OldHash = "c52c7a90e94c5465e55e2c08fa133e27"; NewHash = "7788fa9a9646e7159463bd9b2733690f" }
    @{ Row = 483; Id = "05-050205A"; OldHash = "c13f2fbb39889312ebe51655c8906f50"; NewHash = "4ca5719d81a84971846dd7d5fa957dfe" }
    @{ Row = 485; Id = "05-050314A"; OldHash = "466dbbdd75b29f19f4337f486efc3f50"; NewHash = "11393237a2af2e1122fb08c6fb5ccc47" }
    @{ Row = 506; Id = "05-050202A"; OldHash = "c4e086901e87a390d81c08e4bb9fdebd"; NewHash = "ff555ac7d4b78402ea8ad271f2a4ad9f" }
    @{ Row = 507; Id = "05-050311A"; OldHash = "0609d0fedd5ea75cbda7ed438647537d"; NewHash = "72817237875952ffa7f42e91e248b248" }
    @{ Row = 508; Id = "05-050208TP"; OldHash = "e5689301a7dcef202aae3ff556c77d8e"; NewHash = "6a55751d6462bd11b65b7440271838e8" }
    @{ Row = 513; Id = "05-050306TP"; OldHash = "c3a5a348978c43abef5433eb005a5bcf"; NewHash = "2146a80c226863453b7dfbf1d5243cb5" }
    @{ Row = 521; Id = "05-050317TC"; OldHash = "ddb1b7c47afac42d25035f85f87784b0"; NewHash = "6cd2a825525456c7dfcea9d189a72925" }
    @{ Row = 524; Id = "05-050203A"; OldHash = "51afb509ed48bf7dc9322527ac61fcb8"; NewHash = "477fa6c50530ce476500120b169a8cd4" }
    @{ Row = 532; Id = "05-050317TP"; OldHash = "f233a859f5bd79f16e9f2757fc8bd8f3"; NewHash = "8793a3446d9ad8bc87c6d6c5f6f4fafc" }
    @{ Row = 555; Id = "05-050201A"; OldHash = "56d4c587a2471afc2171571f7d27f0c7"; NewHash = "d801b8e81876e7c4a64433dfd4dc2b7b" }
    @{ Row = 580; Id = "05-050308TP"; OldHash = "c3f60ea1fa19ab1c30e5690afe2c4a50"; NewHash = "7dae245ad2561b7be11a75a20fafd344" }
    @{ Row = 624; Id = "05-050204TP"; OldHash = "49317de9592d0ba2745f2811467e0469"; NewHash = "b84f839796b066a1243c52707d92a25c" }
    @{ Row = 635; Id = "05-050204TC"; OldHash = "ba40ada3c09dc5fca60bdcd71f8eb628"; NewHash = "8cba11ef79a7be606a3e6234747589ea" }
    @{ Row = 637; Id = "05-050302TP"; OldHash = "d800fd7c62ec6c0c850576a265f8b098"; NewHash = "9f944a039a1ab7b1ed0144f8cb9a8842" }
    @{ Row = 657; Id = "05-050313TP"; OldHash = "54a38eade1ac1d27fbf81691c24ba5cb"; NewHash = "e79f1b063b4c5551d8555d5015a20f56" }
    @{ Row = 663; Id = "05-050313TC"; OldHash = "a5cbd6b59e2b8d03a54f6e0298e7be71"; NewHash = "0504657e3c292108bc45a1f21524a881" }
    @{ Row = 674; Id = "05-050317A"; OldHash = "ad397fcbf26ca599eacf01feded20657"; NewHash = "ede95055e968bde6e398118600ce9fc8" }
    @{ Row = 688; Id = "05-050206TP"; OldHash = "be0392d0a6d60ca7e9618cc528ca05c4"; NewHash = "42a5876f6bc762432206df3560f455c8" }
    @{ Row = 693; Id = "05-050206TC"; OldHash = "0f57cd6fedeba799532b3b0ba4b4f37e"; NewHash = "6e28dc225902452ced32c83fd233c074" }
    @{ Row = 708; Id = "05-050304TC"; OldHash = "f50477f816aa871fa3ad5dcebd3db88f"; NewHash = "ecf5019dbf4165b4c5a7e8747ac8d10a" }
    @{ Row = 711; Id = "05-050206A"; OldHash = "a0988d882b880d8a537e3239a8df8a88"; NewHash = "038ced26d73fd13707ae1bd943a110ff" }
    @{ Row = 712; Id = "05-050315A"; OldHash = "b384debc4d274cd463c47814bb03584f"; NewHash = "6b11538b912d342511b93a8d6f30106c" }
    @{ Row = 723; Id = "05-050304TP"; OldHash = "7068eaeabb596cde9800331635f8126a"; NewHash = "d5c064fb0d3303bd2137a09526d79d50" }
    @{ Row = 737; Id = "05-050316A"; OldHash = "facd5956c8a905ffec71ab1f34091b5f"; NewHash = "9be0b7bf49fc4d2b2aea261e627c93ba" }
    @{ Row = 741; Id = "05-050207A"; OldHash = "46fd746591fd901150cf8faf3cfdf069"; NewHash = "934407cd76a3ce7fb54dc59ffc614af9" }
    @{ Row = 750; Id = "05-050315TP"; OldHash = "a8bb0f961ff788315d8ff7476ab13e9e"; NewHash = "f5d456ee9d750b9f00bef309527541dd" }
    @{ Row = 827; Id = "05-050202TP"; OldHash = "8984ed957ef45588ab2b7e250414079d"; NewHash = "67d19e917d77495a0bf1429327bc46f2" }
    @{ Row = 838; Id = "05-050311TC"; OldHash = "8e6ce109c0d4d26ba84457b761a405fb"; NewHash = "9ee254b94c8b6f563f234807c6dd6c15" }
    @{ Row = 843; Id = "05-050311TP"; OldHash = "f8983d626543a44eea837a380749594b"; NewHash = "f52da1897a366f70a7ccd43c07cf406e" }
    @{ Row = 862; Id = "05-050309TC"; OldHash = "489eb9999e6bb05525b849268980c21f"; NewHash = "52df82f178bc48a6dbc0d978371bf77d" }
    @{ Row = 882; Id = "03-030032A"; OldHash = "c9c849f03081bb7a17b5eba5feebb7ea"; NewHash = "d878f735a89572d2273c1e98708e28dd" }
)

foreach ($u in $updates) {
    $idCell = $ws.Cells.Item($u.Row, 1)
    $hashCell = $ws.Cells.Item($u.Row, 2)
    if ($idCell.Value2 -eq $u.Id -and $hashCell.Value2 -eq $u.OldHash) {
        $hashCell.Value = $u.NewHash
    }
}
